# "readworkboox and test 3"
# Swap the header labels in A1/B1 ("name" <-> "sampleid") so the header row
# matches the data actually stored in those columns, nudge the formatting
# of the header cells (A1, B1, D1 - not C1), resize row 1, and move the
# active selection to C1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Swap A1 and B1 header text -------------------------------------------------
$a1Value = $ws.Cells.Item(1, 1).Value2
$b1Value = $ws.Cells.Item(1, 2).Value2

$ws.Cells.Item(1, 1).Value = $b1Value
$ws.Cells.Item(1, 2).Value = $a1Value

# --- Re-apply formatting on the header cells (A1, B1, D1) ----------------------
$ws.Range("A1").Style = "Normal"
$ws.Range("B1").Style = "Normal"
$ws.Range("D1").Style = "Normal"

# --- Row 1 height -----------------------------------------------------------
$ws.Rows.Item(1).RowHeight = 12.8

# --- Move the selection to C1 ------------------------------------------------
$ws.Range("C1").Select()
